$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 12 with the DC motor BOM entry
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "DC motor 31ZY, 6V, 4000rpm"
$ws.Range("D12").Value = "Amazon, ebay, others"
$ws.Range("E12").Value = "https://amzn.eu/d/0FNWCg7"

# Match wrap-text formatting used by the rest of the table's D/E columns
$ws.Range("D12:E12").WrapText = $true

# Row height consistent with the other single-line rows in the table
$ws.Rows.Item(12).RowHeight = 17

# Update selection to reflect the new active cell
$ws.Range("C12").Select()
